$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: Bus3_import with zero capacity/cost and Bus index 3
$ws.Range("A4").Value = "Bus3_import"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 3

# New column D header
$ws.Range("D1").Value = "Bus"
$ws.Range("D1").Font.Bold = $true

# Fill in Bus index for existing rows
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 2

# Update selection to match the authored state
$ws.Range("D11").Select()
